$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 321.61905
$ws.Range("I32").Value = 221.9
$ws.Range("J32").Value = 412.27274
$ws.Range("K32").Value = 221.9
$ws.Range("L32").Value = 412.27274
$ws.Range("M32").Value = 104.1
$ws.Range("N32").Value = -1064.27274
$ws.Range("H62").Value = 1939.875
$ws.Range("I62").Value = 1788.4286
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 1788.4286
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -1164.4286
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 1939.875
$ws.Range("I65").Value = 1788.4286
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 8942.143
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -5822.143
$ws.Range("N65").Value = -21240
$ws.Range("H93").Value = 76633.664
$ws.Range("J93").Value = 76633.664
$ws.Range("L93").Value = 76633.664
$ws.Range("N93").Value = -81625.664
$ws.Range("H98").Value = 1733.7778
$ws.Range("I98").Value = 1700.5
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 1700.5
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = -202.5
$ws.Range("N98").Value = -4996
$ws.Range("H113").Value = 2558
$ws.Range("I113").Value = 1897.2727
$ws.Range("J113").Value = 3284.8
$ws.Range("K113").Value = 1897.2727
$ws.Range("L113").Value = 3284.8
$ws.Range("M113").Value = 1356.7273
$ws.Range("N113").Value = -9792.799999999999
$ws.Range("H116").Value = 2142.8572
$ws.Range("I116").Value = 2000
$ws.Range("J116").Value = 2200
$ws.Range("K116").Value = 2000
$ws.Range("L116").Value = 2200
$ws.Range("M116").Value = 1442
$ws.Range("N116").Value = -9084
$ws.Range("H122").Value = 1733.7778
$ws.Range("I122").Value = 1700.5
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 5101.5
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -2651.5
$ws.Range("N122").Value = -10900
$ws.Range("H141").Value = 1894.0465
$ws.Range("I141").Value = 1253.5834
$ws.Range("J141").Value = 5187.857
$ws.Range("K141").Value = 3760.7502
$ws.Range("L141").Value = 15563.571
$ws.Range("M141").Value = 1419.2498
$ws.Range("N141").Value = -25923.571

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H54").Value = 29500
$ws.Range("J54").Value = 29500
$ws.Range("L54").Value = 29500
$ws.Range("N54").Value = -31038

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H39").Value = 14702
$ws.Range("J39").Value = 14702
$ws.Range("L39").Value = 14702
$ws.Range("N39").Value = -15480

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4865.88
$ws.Range("I16").Value = 6141.5
$ws.Range("J16").Value = 1585.7142
$ws.Range("K16").Value = 6141.5
$ws.Range("L16").Value = 1585.7142
$ws.Range("M16").Value = -5854.5
$ws.Range("N16").Value = -2159.7142
$ws.Range("H105").Value = 970
$ws.Range("I105").Value = 466.66666
$ws.Range("J105").Value = 1725
$ws.Range("K105").Value = 466.66666
$ws.Range("L105").Value = 1725
$ws.Range("M105").Value = 1280.33334
$ws.Range("N105").Value = -5219
$ws.Range("H107").Value = 1006.6
$ws.Range("I107").Value = 1619
$ws.Range("J107").Value = 718.41174
$ws.Range("K107").Value = 1619
$ws.Range("L107").Value = 718.41174
$ws.Range("M107").Value = 301
$ws.Range("N107").Value = -4558.41174
$ws.Range("H113").Value = 4865.88
$ws.Range("I113").Value = 6141.5
$ws.Range("J113").Value = 1585.7142
$ws.Range("K113").Value = 6141.5
$ws.Range("L113").Value = 1585.7142
$ws.Range("M113").Value = -3971.5
$ws.Range("N113").Value = -5925.7142
$ws.Range("H117").Value = 56156
$ws.Range("J117").Value = 56156
$ws.Range("L117").Value = 56156
$ws.Range("N117").Value = -65334
$ws.Range("H132").Value = 1701.4546
$ws.Range("I132").Value = 1169.8684
$ws.Range("J132").Value = 2889.7058
$ws.Range("K132").Value = 3509.6052
$ws.Range("L132").Value = 8669.117400000001
$ws.Range("M132").Value = -979.6052
$ws.Range("N132").Value = -13729.1174
$ws.Range("H134").Value = 2814.9
$ws.Range("I134").Value = 1611.8928
$ws.Range("J134").Value = 4346
$ws.Range("K134").Value = 4835.678400000001
$ws.Range("L134").Value = 13038
$ws.Range("M134").Value = -2300.678400000001
$ws.Range("N134").Value = -18108

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 9848716
$ws.Range("I9").Value = 2000
$ws.Range("J9").Value = 10464136
$ws.Range("K9").Value = 6000
$ws.Range("L9").Value = 31392408
$ws.Range("M9").Value = -5776
$ws.Range("N9").Value = -31392856
$ws.Range("H122").Value = 862.6585
$ws.Range("I122").Value = 404.85715
$ws.Range("J122").Value = 956.91174
$ws.Range("K122").Value = 3643.71435
$ws.Range("L122").Value = 8612.20566
$ws.Range("M122").Value = -1193.71435
$ws.Range("N122").Value = -13512.20566
$ws.Range("H131").Value = 34343.93
$ws.Range("I131").Value = 1644.2307
$ws.Range("J131").Value = 60912.438
$ws.Range("K131").Value = 4932.6921
$ws.Range("L131").Value = 182737.314
$ws.Range("M131").Value = 107.3078999999998
$ws.Range("N131").Value = -192817.314

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 49999.8
$ws.Range("J104").Value = 49999.8
$ws.Range("L104").Value = 49999.8
$ws.Range("N104").Value = -56987.8
$ws.Range("H112").Value = 49900
$ws.Range("J112").Value = 49900
$ws.Range("L112").Value = 49900
$ws.Range("N112").Value = -52116

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 37000
$ws.Range("H3").Value = 12500
$ws.Range("I3").Value = 10000
$ws.Range("J3").Value = 15000
$ws.Range("K3").Value = 10000
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = -9888
$ws.Range("N3").Value = -15224
$ws.Range("H9").Value = 3333.3333
$ws.Range("I9").Value = 500
$ws.Range("J9").Value = 9000
$ws.Range("K9").Value = 500
$ws.Range("L9").Value = 9000
$ws.Range("M9").Value = -276
$ws.Range("N9").Value = -9448
$ws.Range("H10").Value = 4483.8335
$ws.Range("I10").Value = 1034.3334
$ws.Range("J10").Value = 7933.3335
$ws.Range("K10").Value = 1034.3334
$ws.Range("L10").Value = 7933.3335
$ws.Range("M10").Value = -894.3334
$ws.Range("N10").Value = -8213.333500000001
$ws.Range("H14").Value = 7000
$ws.Range("I14").Value = 1000
$ws.Range("K14").Value = 1000
$ws.Range("M14").Value = -828
$ws.Range("H15").Value = 12500
$ws.Range("I15").Value = 10000
$ws.Range("J15").Value = 15000
$ws.Range("K15").Value = 10000
$ws.Range("L15").Value = 15000
$ws.Range("M15").Value = -9830
$ws.Range("N15").Value = -15340
$ws.Range("H21").Value = 5024.4165
$ws.Range("J21").Value = 15000
$ws.Range("L21").Value = 15000
$ws.Range("H22").Value = 856.5
$ws.Range("J22").Value = 842.2308
$ws.Range("L22").Value = 842.2308
$ws.Range("N22").Value = -1432.2308
$ws.Range("H27").Value = 856.5
$ws.Range("J27").Value = 842.2308
$ws.Range("L27").Value = 842.2308
$ws.Range("N27").Value = -1056.2308
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H51").Value = 20084
$ws.Range("J51").Value = 20084
$ws.Range("L51").Value = 20084
$ws.Range("N51").Value = -21040
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H61").Value = 353896.8
$ws.Range("I61").Value = 12294.4
$ws.Range("J61").Value = 1113013.4
$ws.Range("K61").Value = 12294.4
$ws.Range("L61").Value = 1113013.4
$ws.Range("M61").Value = -12092.4
$ws.Range("N61").Value = -1113417.4
$ws.Range("H93").Value = 1082.5294
$ws.Range("I93").Value = 876.1539
$ws.Range("J93").Value = 1753.25
$ws.Range("K93").Value = 876.1539
$ws.Range("L93").Value = 1753.25
$ws.Range("M93").Value = 371.8461
$ws.Range("N93").Value = -4249.25
$ws.Range("H111").Value = 39800
$ws.Range("J111").Value = 39800
$ws.Range("L111").Value = 39800
$ws.Range("N111").Value = -47980
$ws.Range("H113").Value = 353896.8
$ws.Range("I113").Value = 12294.4
$ws.Range("J113").Value = 1113013.4
$ws.Range("K113").Value = 12294.4
$ws.Range("L113").Value = 1113013.4
$ws.Range("M113").Value = -10124.4
$ws.Range("N113").Value = -1117353.4

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 16400.6
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 16400.6
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 16400.6
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -16624.6
$ws.Range("H4").Value = 2745.4546
$ws.Range("I4").Value = 200
$ws.Range("K4").Value = 200
$ws.Range("H41").Value = 12728.111
$ws.Range("J41").Value = 12728.111
$ws.Range("L41").Value = 12728.111
$ws.Range("N41").Value = -13508.111
$ws.Range("H100").Value = 1074.8823
$ws.Range("I100").Value = 496.66666
$ws.Range("K100").Value = 993.33332
$ws.Range("M100").Value = -452.33332
$ws.Range("H113").Value = 4126.6665
$ws.Range("I113").Value = 6700.4375
$ws.Range("J113").Value = 1185.2142
$ws.Range("K113").Value = 20101.3125
$ws.Range("L113").Value = 3555.6426
$ws.Range("M113").Value = -17931.3125
$ws.Range("N113").Value = -7895.642599999999
$ws.Range("H126").Value = 1446.3846
$ws.Range("I126").Value = 1463.4546
$ws.Range("J126").Value = 1352.5
$ws.Range("K126").Value = 4390.3638
$ws.Range("L126").Value = 4057.5
$ws.Range("M126").Value = -1920.3638
$ws.Range("N126").Value = -8997.5
